# Auto-generated edit script applying the diff changes to before.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 64.5
$ws.Range("J6").Value = 111
$ws.Range("K6").Value = 193.5
$ws.Range("L6").Value = 333
$ws.Range("M6").Value = -81.5
$ws.Range("N6").Value = -557
$ws.Range("H17").Value = 658.4375
$ws.Range("J17").Value = 658.29034
$ws.Range("L17").Value = 1974.87102
$ws.Range("N17").Value = -2310.87102
$ws.Range("H33").Value = 356.42856
$ws.Range("I33").Value = 309.73685
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 309.73685
$ws.Range("L33").Value = 800
$ws.Range("M33").Value = -80.73685
$ws.Range("N33").Value = -1258
$ws.Range("H105").Value = 40402.2
$ws.Range("J105").Value = 40402.2
$ws.Range("L105").Value = 40402.2
$ws.Range("N105").Value = -47390.2
$ws.Range("H129").Value = 887.02856
$ws.Range("J129").Value = 898.6866
$ws.Range("L129").Value = 2696.0598
$ws.Range("N129").Value = -12696.0598
$ws.Range("H138").Value = 3068.743
$ws.Range("J138").Value = 3359.6667
$ws.Range("L138").Value = 10079.0001
$ws.Range("N138").Value = -20359.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2198.2031
$ws.Range("I32").Value = 1464.0176
$ws.Range("J32").Value = 8176.5713
$ws.Range("K32").Value = 1464.0176
$ws.Range("L32").Value = 8176.5713
$ws.Range("M32").Value = -1177.0176
$ws.Range("N32").Value = -8750.5713
$ws.Range("H61").Value = 3665.2942
$ws.Range("I61").Value = 2630.2856
$ws.Range("J61").Value = 4389.8
$ws.Range("K61").Value = 2630.2856
$ws.Range("L61").Value = 4389.8
$ws.Range("M61").Value = -2418.2856
$ws.Range("N61").Value = -4813.8
$ws.Range("H104").Value = 39000
$ws.Range("J104").Value = 39000
$ws.Range("L104").Value = 39000
$ws.Range("N104").Value = -45988
$ws.Range("H122").Value = 1736.84
$ws.Range("I122").Value = 1656.75
$ws.Range("J122").Value = 2057.2
$ws.Range("K122").Value = 4970.25
$ws.Range("L122").Value = 6171.599999999999
$ws.Range("M122").Value = -2520.25
$ws.Range("N122").Value = -11071.6
$ws.Range("H136").Value = 3665.2942
$ws.Range("I136").Value = 2630.2856
$ws.Range("J136").Value = 4389.8
$ws.Range("K136").Value = 7890.8568
$ws.Range("L136").Value = 13169.4
$ws.Range("M136").Value = -5340.8568
$ws.Range("N136").Value = -18269.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3113.7585
$ws.Range("I134").Value = 3557.7083
$ws.Range("J134").Value = 982.8
$ws.Range("K134").Value = 10673.1249
$ws.Range("L134").Value = 2948.4
$ws.Range("M134").Value = -8138.124899999999
$ws.Range("N134").Value = -8018.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 35076
$ws.Range("J43").Value = 35076
$ws.Range("L43").Value = 35076
$ws.Range("N43").Value = -35444
$ws.Range("H94").Value = 5969.5
$ws.Range("I94").Value = 2875
$ws.Range("J94").Value = 9064
$ws.Range("K94").Value = 2875
$ws.Range("L94").Value = 9064
$ws.Range("M94").Value = -2424
$ws.Range("N94").Value = -9966
$ws.Range("H101").Value = 35076
$ws.Range("J101").Value = 35076
$ws.Range("L101").Value = 35076
$ws.Range("N101").Value = -41566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 99.5
$ws.Range("J7").Value = 99.5
$ws.Range("L7").Value = 298.5
$ws.Range("N7").Value = -522.5
$ws.Range("H80").Value = 1750
$ws.Range("I80").Value = 1750
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5250
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -4314
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1750
$ws.Range("I83").Value = 1750
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15750
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11070
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 899
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 899
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H92").Value = 664.44446
$ws.Range("J92").Value = 622.5
$ws.Range("L92").Value = 1867.5
$ws.Range("N92").Value = -4363.5
$ws.Range("H131").Value = 802.86
$ws.Range("J131").Value = 823.1042
$ws.Range("L131").Value = 2469.3126
$ws.Range("N131").Value = -12549.3126
$ws.Range("H137").Value = 2966.5293
$ws.Range("I137").Value = 500
$ws.Range("J137").Value = 3295.4
$ws.Range("K137").Value = 1500
$ws.Range("L137").Value = 9886.200000000001
$ws.Range("M137").Value = 3600
$ws.Range("N137").Value = -20086.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 27996
$ws.Range("J39").Value = 27996
$ws.Range("L39").Value = 27996
$ws.Range("N39").Value = -29060
$ws.Range("H101").Value = 31666.334
$ws.Range("J101").Value = 31666.334
$ws.Range("L101").Value = 31666.334
$ws.Range("N101").Value = -38156.334
$ws.Range("H126").Value = 4806.552
$ws.Range("I126").Value = 4111.7646
$ws.Range("J126").Value = 5790.8335
$ws.Range("K126").Value = 12335.2938
$ws.Range("L126").Value = 17372.5005
$ws.Range("M126").Value = -9865.293800000001
$ws.Range("N126").Value = -22312.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H98").Value = 48266.668
$ws.Range("J98").Value = 48266.668
$ws.Range("L98").Value = 48266.668
$ws.Range("N98").Value = -54256.668
$ws.Range("H104").Value = 31792.5
$ws.Range("J104").Value = 31792.5
$ws.Range("L104").Value = 31792.5
$ws.Range("N104").Value = -38780.5
$ws.Range("H132").Value = 2090.75
$ws.Range("I132").Value = 1527.8
$ws.Range("J132").Value = 3779.6
$ws.Range("K132").Value = 4583.4
$ws.Range("L132").Value = 11338.8
$ws.Range("M132").Value = -2053.4
$ws.Range("N132").Value = -16398.8
